$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) data cells to remain text, matching the
# original inline-string cell type (avoids Excel auto-converting
# numeric-looking strings like "1.000" or "0.9999" into numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.416.54'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.847.83'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '240.72'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = '0.6278'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '0.07683'
$ws.Range("E8").Value = '  +2.12%  '
$ws.Range("D9").Value = '0.2919'
$ws.Range("E9").Value = '  -0.38%  '
$ws.Range("D10").Value = '24.75'
$ws.Range("E10").Value = '  +1.50%  '
$ws.Range("D11").Value = '0.07735'
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = '1.847.43'
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").Value = '5.025'
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("D14").Value = '0.00001079'
$ws.Range("E14").Value = '  +3.55%  '
$ws.Range("D15").Value = '0.6795'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").Value = '83.48'
$ws.Range("E16").Value = '  +0.79%  '
$ws.Range("D17").Value = '6.171'
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").Value = '29.442.92'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").Value = '227.75'
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").Value = '12.41'
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '7.404'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '157.61'
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("D25").Value = '0.1376'
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("D26").Value = '8.391'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").Value = '17.67'
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").Value = '1.343'
$ws.Range("E28").Value = '  +5.26%  '
$ws.Range("D29").Value = '1.465'
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("D30").Value = '0.05689'
$ws.Range("E30").Value = '  +1.01%  '
$ws.Range("D31").Value = '4.115'
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("D32").Value = '4.027'
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("D33").Value = '1.841'
$ws.Range("E33").Value = '  +0.40%  '
$ws.Range("D34").Value = '1.161'
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("D35").Value = '0.7078'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").Value = '2.584'
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").Value = '2.777'
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("D38").Value = '1.226.45'
$ws.Range("E38").Value = '  -1.02%  '
$ws.Range("D39").Value = '0.01790'
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("D40").Value = '6.534'
$ws.Range("E40").Value = '  +4.64%  '
$ws.Range("D41").Value = '0.9051'
$ws.Range("E41").Value = '  +0.71%  '
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").Value = '101.68'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").Value = '66.09'
$ws.Range("E44").Value = '  +1.10%  '
$ws.Range("D45").Value = '0.00000000121'
$ws.Range("E45").Value = '  +1.68%  '
$ws.Range("D46").Value = '7.150'
$ws.Range("E46").Value = '  +1.93%  '
$ws.Range("D47").Value = '0.4014'
$ws.Range("E47").Value = '  +0.62%  '
$ws.Range("D48").Value = '8.975'
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("D49").Value = '0.1145'
$ws.Range("E49").Value = '  +2.28%  '
$ws.Range("D50").Value = '1.673'
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("D51").Value = '0.05714'
$ws.Range("E51").Value = '  +0.07%  '
